# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The underlying change: the "Valor Mora" figures for the first and last
# period rows of the account-statement table were corrected/updated -
# F16 (period 2410) and F22 (period 2404) had their values swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 15600
